# Edit script to apply the "Horarios actualizados Linea 141 - 983" update
# Updates 3 worksheets: LP1912, LP1912-215, 6203-6173
# Each sheet has:
#   Row1: title, Row2: "Ultima actualizacion: HH:MM:SS", Row3: "Total filas: N"
#   Row5: column headers, Row6+: data rows (Hora_Scrap, Hora_Llegada, Linea, Minutos, Parada)

$wb = $excel.ActiveWorkbook

$nuevaHora = "15:59:49"

# ---------------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Cells.Item(2,1).Value = "Última actualización: $nuevaHora"
$ws1.Cells.Item(3,1).Value = "Total filas: 260"

$data1 = @(
    @("05:44:02","05:44","14_ABASTO",0,"LP1912"),
    @("05:44:02","05:47","17_ROMERO",3,"LP1912"),
    @("05:58:04","06:09","10_OLMOS",11,"LP1912"),
    @("05:58:04","06:16","215A_EL PATO",18,"LP1912"),
    @("06:19:59","06:29","23_HERNANDEZ",10,"LP1912"),
    @("05:58:04","06:30","23_HERNANDEZ",32,"LP1912"),
    @("06:19:59","06:33","11_ETCHEVERRY",14,"LP1912"),
    @("05:58:04","06:34","11_ETCHEVERRY",36,"LP1912"),
    @("06:38:54","06:38","17X38_ROMERO",0,"LP1912"),
    @("05:58:04","06:39","17X38_ROMERO",41,"LP1912"),
    @("05:44:02","06:40","17X38_ROMERO",56,"LP1912"),
    @("06:38:54","06:40","16_SANTA ANA",2,"LP1912"),
    @("05:58:04","06:41","16_SANTA ANA",43,"LP1912"),
    @("06:38:54","06:56","215A_EL PATO",18,"LP1912"),
    @("06:56:24","06:57","215A_EL PATO",1,"LP1912"),
    @("06:38:54","06:58","225_GOMEZ",20,"LP1912"),
    @("06:56:24","06:59","225_GOMEZ",3,"LP1912"),
    @("07:15:48","07:15","215C_EL PATO",0,"LP1912"),
    @("06:56:24","07:16","215C_EL PATO",20,"LP1912"),
    @("06:38:54","07:18","14_ABASTO",40,"LP1912"),
    @("07:15:48","07:19","14_ABASTO",4,"LP1912"),
    @("07:15:48","07:20","16_SANTA ANA",5,"LP1912"),
    @("06:56:24","07:21","16_SANTA ANA",25,"LP1912"),
    @("07:15:48","07:21","23_HERNANDEZ",6,"LP1912"),
    @("06:56:24","07:22","23_HERNANDEZ",26,"LP1912"),
    @("07:15:48","07:29","17X38_ROMERO",14,"LP1912"),
    @("06:38:54","07:34","10_OLMOS",56,"LP1912"),
    @("07:15:48","07:35","10_OLMOS",20,"LP1912"),
    @("07:15:48","07:36","27_EL RETIRO",21,"LP1912"),
    @("06:56:24","07:37","27_EL RETIRO",41,"LP1912"),
    @("07:15:48","07:43","215A_EL PATO",28,"LP1912"),
    @("07:40:11","07:44","215A_EL PATO",4,"LP1912"),
    @("07:52:32","07:54","14_ABASTO",2,"LP1912"),
    @("07:40:11","07:55","14_ABASTO",15,"LP1912"),
    @("06:19:59","07:59","17_ROMERO",100,"LP1912"),
    @("06:38:54","08:00","16_SANTA ANA",82,"LP1912"),
    @("07:52:32","08:00","17_ROMERO",8,"LP1912"),
    @("07:52:32","08:01","16_SANTA ANA",9,"LP1912"),
    @("07:52:32","08:06","23_HERNANDEZ",14,"LP1912"),
    @("08:10:22","08:11","10_OLMOS",1,"LP1912"),
    @("06:19:59","08:12","15X38_ABASTO",113,"LP1912"),
    @("08:10:22","08:13","15X38_ABASTO",3,"LP1912"),
    @("06:38:54","08:28","11_ETCHEVERRY",110,"LP1912"),
    @("08:10:22","08:29","15_ABASTO",19,"LP1912"),
    @("08:10:22","08:29","11_ETCHEVERRY",19,"LP1912"),
    @("08:30:14","08:30","11_ETCHEVERRY",0,"LP1912"),
    @("08:30:14","08:30","15_ABASTO",0,"LP1912"),
    @("08:40:59","08:40","16_P MOR-SANTA ANA",0,"LP1912"),
    @("08:40:59","08:40","10_OLMOS",0,"LP1912"),
    @("07:52:32","08:41","10_OLMOS",49,"LP1912"),
    @("08:30:14","08:41","16_P MOR-SANTA ANA",11,"LP1912"),
    @("08:40:59","08:43","215C_EL PATO",3,"LP1912"),
    @("08:30:14","08:44","215C_EL PATO",14,"LP1912"),
    @("07:40:11","08:46","23_HERNANDEZ",66,"LP1912"),
    @("07:52:32","08:47","23_HERNANDEZ",55,"LP1912"),
    @("08:40:59","08:51","23_HERNANDEZ",11,"LP1912"),
    @("08:52:33","08:52","23_HERNANDEZ",0,"LP1912"),
    @("08:52:33","08:53","215B_EL PATO",1,"LP1912"),
    @("08:30:14","08:54","23_HERNANDEZ",24,"LP1912"),
    @("08:30:14","08:54","215B_EL PATO",24,"LP1912"),
    @("08:40:59","08:57","215A_EL PATO",17,"LP1912"),
    @("08:52:33","08:58","215A_EL PATO",6,"LP1912"),
    @("08:40:59","09:04","10_OLMOS",24,"LP1912"),
    @("08:52:33","09:05","10_OLMOS",13,"LP1912"),
    @("08:52:33","09:06","16_SANTA ANA",14,"LP1912"),
    @("08:10:22","09:11","27_EL RETIRO",61,"LP1912"),
    @("07:52:32","09:12","27_EL RETIRO",80,"LP1912"),
    @("07:40:11","09:14","27_EL RETIRO",94,"LP1912"),
    @("08:52:33","09:16","27_EL RETIRO",24,"LP1912"),
    @("08:30:14","09:17","27_EL RETIRO",47,"LP1912"),
    @("08:40:59","09:17","15X38_ABASTO",37,"LP1912"),
    @("08:52:33","09:17","14_ABASTO",25,"LP1912"),
    @("08:52:33","09:18","15X38_ABASTO",26,"LP1912"),
    @("08:30:14","09:18","14_ABASTO",48,"LP1912"),
    @("09:23:52","09:31","16_SANTA ANA",8,"LP1912"),
    @("09:23:52","09:36","23_HERNANDEZ",13,"LP1912"),
    @("08:40:59","09:38","15_ABASTO",58,"LP1912"),
    @("09:23:52","09:39","15_ABASTO",16,"LP1912"),
    @("09:23:52","09:41","11_ETCHEVERRY",18,"LP1912"),
    @("08:30:14","09:42","11_ETCHEVERRY",72,"LP1912"),
    @("08:40:59","09:42","16_P MOR-SANTA ANA",62,"LP1912"),
    @("09:23:52","09:43","16_P MOR-SANTA ANA",20,"LP1912"),
    @("08:40:59","09:52","10_OLMOS",72,"LP1912"),
    @("09:23:52","09:53","10_OLMOS",30,"LP1912"),
    @("09:23:52","09:58","215C_EL PATO",35,"LP1912"),
    @("08:52:33","09:59","215C_EL PATO",67,"LP1912"),
    @("09:23:52","10:05","14_ABASTO",42,"LP1912"),
    @("08:30:14","10:06","14_ABASTO",96,"LP1912"),
    @("08:40:59","10:12","17X38_ROMERO",92,"LP1912"),
    @("10:07:51","10:13","17X38_ROMERO",6,"LP1912"),
    @("09:23:52","10:15","23_HERNANDEZ",52,"LP1912"),
    @("08:52:33","10:21","23_HERNANDEZ",89,"LP1912"),
    @("10:07:51","10:22","23_HERNANDEZ",15,"LP1912"),
    @("08:30:14","10:24","23_HERNANDEZ",114,"LP1912"),
    @("10:07:51","10:25","16_SANTA ANA",18,"LP1912"),
    @("08:40:59","10:28","15_ABASTO",108,"LP1912"),
    @("10:07:51","10:29","14_ABASTO",22,"LP1912"),
    @("10:07:51","10:29","15_ABASTO",22,"LP1912"),
    @("10:07:51","10:43","11X44_ETCHEVERRY",36,"LP1912"),
    @("10:41:48","10:44","11X44_ETCHEVERRY",3,"LP1912"),
    @("10:41:48","10:46","15_P INDUSTRIAL",5,"LP1912"),
    @("09:23:52","10:53","27_EL RETIRO",90,"LP1912"),
    @("10:41:48","10:55","16_SANTA ANA",14,"LP1912"),
    @("10:41:48","10:56","27_EL RETIRO",15,"LP1912"),
    @("10:56:01","10:57","27_EL RETIRO",1,"LP1912"),
    @("10:56:01","10:59","10_OLMOS",3,"LP1912"),
    @("10:56:01","11:01","81_EL PELIGRO",5,"LP1912"),
    @("09:23:52","11:01","10_OLMOS",98,"LP1912"),
    @("10:07:51","11:03","23_HERNANDEZ",56,"LP1912"),
    @("10:56:01","11:07","23_HERNANDEZ",11,"LP1912"),
    @("10:56:01","11:10","16_P MOR-SANTA ANA",14,"LP1912"),
    @("10:41:48","11:14","14_ABASTO",33,"LP1912"),
    @("11:13:01","11:15","14_ABASTO",2,"LP1912"),
    @("11:13:01","11:15","15X38_ABASTO",2,"LP1912"),
    @("09:23:52","11:19","10_OLMOS",116,"LP1912"),
    @("10:07:51","11:21","16_SANTA ANA",74,"LP1912"),
    @("11:13:01","11:25","16_SANTA ANA",12,"LP1912"),
    @("11:13:01","11:29","10_OLMOS",16,"LP1912"),
    @("10:41:48","11:30","215C_EL PATO",49,"LP1912"),
    @("11:13:01","11:31","215C_EL PATO",18,"LP1912"),
    @("11:35:40","11:41","215B_EL PATO",6,"LP1912"),
    @("11:13:01","11:42","215B_EL PATO",29,"LP1912"),
    @("11:35:40","11:45","15X38_ABASTO",10,"LP1912"),
    @("10:56:01","11:45","23_HERNANDEZ",49,"LP1912"),
    @("11:48:20","11:48","16_SANTA ANA",0,"LP1912"),
    @("10:41:48","11:49","23_HERNANDEZ",68,"LP1912"),
    @("10:07:51","11:52","225_GOMEZ",105,"LP1912"),
    @("11:48:20","11:52","23_HERNANDEZ",4,"LP1912"),
    @("11:48:20","11:53","225_GOMEZ",5,"LP1912"),
    @("11:35:40","11:53","23_HERNANDEZ",18,"LP1912"),
    @("11:55:01","11:55","225_GOMEZ",0,"LP1912"),
    @("11:55:01","11:58","17_ROMERO",3,"LP1912"),
    @("11:35:40","12:05","11_ETCHEVERRY",30,"LP1912"),
    @("11:13:01","12:06","11_ETCHEVERRY",53,"LP1912"),
    @("11:55:01","12:10","16_P MOR-SANTA ANA",15,"LP1912"),
    @("11:55:01","12:10","15_ABASTO",15,"LP1912"),
    @("12:12:04","12:12","16_SANTA ANA",0,"LP1912"),
    @("12:12:04","12:13","15_ABASTO",1,"LP1912"),
    @("12:12:04","12:17","10_OLMOS",5,"LP1912"),
    @("11:35:40","12:21","215C_EL PATO",46,"LP1912"),
    @("12:12:04","12:22","215C_EL PATO",10,"LP1912"),
    @("11:35:40","12:31","23_HERNANDEZ",56,"LP1912"),
    @("11:13:01","12:31","27_EL RETIRO",78,"LP1912"),
    @("11:48:20","12:32","23_HERNANDEZ",44,"LP1912"),
    @("11:48:20","12:32","14_ABASTO",44,"LP1912"),
    @("11:13:01","12:33","23_HERNANDEZ",80,"LP1912"),
    @("12:33:54","12:33","14_ABASTO",0,"LP1912"),
    @("10:56:01","12:34","27_EL RETIRO",98,"LP1912"),
    @("12:33:54","12:34","15_ABASTO",1,"LP1912"),
    @("12:33:54","12:36","23_HERNANDEZ",3,"LP1912"),
    @("12:33:54","12:36","27_EL RETIRO",3,"LP1912"),
    @("12:12:04","12:37","23_HERNANDEZ",25,"LP1912"),
    @("12:12:04","12:37","27_EL RETIRO",25,"LP1912"),
    @("12:47:00","12:47","14_ABASTO",0,"LP1912"),
    @("12:33:54","12:48","16_SANTA ANA",15,"LP1912"),
    @("11:55:01","12:48","14_ABASTO",53,"LP1912"),
    @("12:47:00","12:48","15X38_ABASTO",1,"LP1912"),
    @("12:54:41","13:02","11_ETCHEVERRY",8,"LP1912"),
    @("12:54:41","13:03","215C_EL PATO",9,"LP1912"),
    @("12:47:00","13:03","11_ETCHEVERRY",16,"LP1912"),
    @("12:47:00","13:04","215C_EL PATO",17,"LP1912"),
    @("12:54:41","13:13","16_SANTA ANA",19,"LP1912"),
    @("13:14:41","13:17","10_OLMOS",3,"LP1912"),
    @("13:14:41","13:19","15_ABASTO",5,"LP1912"),
    @("13:14:41","13:21","23_HERNANDEZ",7,"LP1912"),
    @("12:47:00","13:22","23_HERNANDEZ",35,"LP1912"),
    @("12:12:04","13:24","23_HERNANDEZ",72,"LP1912"),
    @("13:14:41","13:25","16_P MOR-SANTA ANA",11,"LP1912"),
    @("13:14:41","13:32","14_ABASTO",18,"LP1912"),
    @("12:33:54","13:32","215A_EL PATO",59,"LP1912"),
    @("12:47:00","13:33","14_ABASTO",46,"LP1912"),
    @("13:14:41","13:33","215A_EL PATO",19,"LP1912"),
    @("13:43:25","13:43","225_GOMEZ",0,"LP1912"),
    @("12:33:54","13:46","225_GOMEZ",73,"LP1912"),
    @("13:14:41","13:47","225_GOMEZ",33,"LP1912"),
    @("11:55:01","13:49","11_ETCHEVERRY",114,"LP1912"),
    @("13:43:25","13:54","15_ABASTO",11,"LP1912"),
    @("12:54:41","13:57","16_SANTA ANA",63,"LP1912"),
    @("12:47:00","13:58","16_SANTA ANA",71,"LP1912"),
    @("12:33:54","14:01","10_OLMOS",88,"LP1912"),
    @("13:14:41","14:02","23_HERNANDEZ",48,"LP1912"),
    @("13:43:25","14:02","16_SANTA ANA",19,"LP1912"),
    @("13:57:31","14:02","10_OLMOS",5,"LP1912"),
    @("13:57:31","14:06","23_HERNANDEZ",9,"LP1912"),
    @("12:47:00","14:07","23_HERNANDEZ",80,"LP1912"),
    @("13:57:31","14:14","15_ABASTO",17,"LP1912"),
    @("14:15:21","14:15","15_ABASTO",0,"LP1912"),
    @("13:57:31","14:16","27_EL RETIRO",19,"LP1912"),
    @("14:15:21","14:17","27_EL RETIRO",2,"LP1912"),
    @("13:57:31","14:17","11_ETCHEVERRY",20,"LP1912"),
    @("12:47:00","14:18","11_ETCHEVERRY",91,"LP1912"),
    @("14:15:21","14:19","11_ETCHEVERRY",4,"LP1912"),
    @("14:15:21","14:27","16_SANTA ANA",12,"LP1912"),
    @("12:33:54","14:31","14X44_ABASTO",118,"LP1912"),
    @("14:15:21","14:32","14X44_ABASTO",17,"LP1912"),
    @("14:33:43","14:33","215C_EL PATO",0,"LP1912"),
    @("14:33:43","14:33","14X44_ABASTO",0,"LP1912"),
    @("14:15:21","14:34","215C_EL PATO",19,"LP1912"),
    @("14:33:43","14:39","16_P MOR-SANTA ANA",6,"LP1912"),
    @("14:33:43","14:47","215B_EL PATO",14,"LP1912"),
    @("14:47:05","14:48","215B_EL PATO",1,"LP1912"),
    @("14:33:43","14:51","23_HERNANDEZ",18,"LP1912"),
    @("14:47:05","14:51","16_SANTA ANA",4,"LP1912"),
    @("14:47:05","14:52","23_HERNANDEZ",5,"LP1912"),
    @("14:53:58","14:53","215_EL PELIGRO",0,"LP1912"),
    @("14:47:05","14:54","215_EL PELIGRO",7,"LP1912"),
    @("14:33:43","15:02","10_OLMOS",29,"LP1912"),
    @("13:14:41","15:11","14_ABASTO",117,"LP1912"),
    @("14:53:58","15:13","17X38_ROMERO",20,"LP1912"),
    @("14:53:58","15:16","16_SANTA ANA",23,"LP1912"),
    @("15:17:56","15:17","14_ABASTO",0,"LP1912"),
    @("14:47:05","15:18","14_ABASTO",31,"LP1912"),
    @("15:17:56","15:33","215C_EL PATO",16,"LP1912"),
    @("14:47:05","15:34","215C_EL PATO",47,"LP1912"),
    @("15:17:56","15:36","23_HERNANDEZ",19,"LP1912"),
    @("15:17:56","15:41","11_ETCHEVERRY",24,"LP1912"),
    @("14:47:05","15:42","11_ETCHEVERRY",55,"LP1912"),
    @("14:53:58","15:52","16_P MOR-SANTA ANA",59,"LP1912"),
    @("15:17:56","15:52","10_OLMOS",35,"LP1912"),
    @("14:53:58","15:52","15X38_ABASTO",59,"LP1912"),
    @("14:33:43","15:52","27_EL RETIRO",79,"LP1912"),
    @("15:47:47","15:53","15X38_ABASTO",6,"LP1912"),
    @("14:47:05","15:53","27_EL RETIRO",66,"LP1912"),
    @("15:47:47","15:53","16_P MOR-SANTA ANA",6,"LP1912"),
    @("15:47:47","15:53","10_OLMOS",6,"LP1912"),
    @("15:47:47","15:56","27_EL RETIRO",9,"LP1912"),
    @("14:15:21","15:57","27_EL RETIRO",102,"LP1912"),
    @("14:53:58","16:04","16_SANTA ANA",71,"LP1912"),
    @("15:59:49","16:05","14_ABASTO",6,"LP1912"),
    @("15:59:49","16:05","16_SANTA ANA",6,"LP1912"),
    @("14:47:05","16:06","14_ABASTO",79,"LP1912"),
    @("14:53:58","16:13","17_ROMERO",80,"LP1912"),
    @("15:59:49","16:14","17_ROMERO",15,"LP1912"),
    @("15:17:56","16:16","10_OLMOS",59,"LP1912"),
    @("15:59:49","16:17","10_OLMOS",18,"LP1912"),
    @("15:59:49","16:21","23_HERNANDEZ",22,"LP1912"),
    @("15:17:56","16:22","23_HERNANDEZ",65,"LP1912"),
    @("15:59:49","16:29","14_ABASTO",30,"LP1912"),
    @("15:47:47","16:30","16_SANTA ANA",43,"LP1912"),
    @("14:53:58","16:33","83_ALUAR",100,"LP1912"),
    @("15:59:49","16:34","83_ALUAR",35,"LP1912"),
    @("15:17:56","16:40","225_GOMEZ",83,"LP1912"),
    @("15:59:49","16:41","225_GOMEZ",42,"LP1912"),
    @("14:53:58","16:45","17_ROMERO",112,"LP1912"),
    @("15:17:56","16:46","17_ROMERO",89,"LP1912"),
    @("15:59:49","16:53","11_ETCHEVERRY",54,"LP1912"),
    @("15:59:49","16:55","16_SANTA ANA",56,"LP1912"),
    @("15:59:49","16:58","15_ABASTO",59,"LP1912"),
    @("15:59:49","17:06","23_HERNANDEZ",67,"LP1912"),
    @("15:47:47","17:07","23_HERNANDEZ",80,"LP1912"),
    @("15:59:49","17:07","16_P MOR-SANTA ANA",68,"LP1912"),
    @("15:59:49","17:08","10_OLMOS",69,"LP1912"),
    @("15:59:49","17:09","215C_EL PATO",70,"LP1912"),
    @("15:59:49","17:21","15X38_ABASTO",82,"LP1912"),
    @("15:59:49","17:33","27_EL RETIRO",94,"LP1912"),
    @("15:59:49","17:34","17_ROMERO",95,"LP1912"),
    @("15:47:47","17:36","27_EL RETIRO",109,"LP1912"),
    @("15:59:49","17:38","215B_EL PATO",99,"LP1912"),
    @("15:47:47","17:42","215_EL PELIGRO",115,"LP1912"),
    @("15:59:49","17:45","215_EL PELIGRO",106,"LP1912")
)

for ($i = 0; $i -lt $data1.Length; $i++) {
    $row = 6 + $i
    $r = $data1[$i]
    $ws1.Cells.Item($row, 1).Value = $r[0]
    $ws1.Cells.Item($row, 2).Value = $r[1]
    $ws1.Cells.Item($row, 3).Value = $r[2]
    $ws1.Cells.Item($row, 4).Value = $r[3]
    $ws1.Cells.Item($row, 5).Value = $r[4]
}

# ---------------------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Cells.Item(2,1).Value = "Última actualización: $nuevaHora"
$ws2.Cells.Item(3,1).Value = "Total filas: 37"

$data2 = @(
    @("05:58:04","06:16","215A_EL PATO",18,"LP1912"),
    @("06:38:54","06:56","215A_EL PATO",18,"LP1912"),
    @("06:56:24","06:57","215A_EL PATO",1,"LP1912"),
    @("07:15:48","07:15","215C_EL PATO",0,"LP1912"),
    @("06:56:24","07:16","215C_EL PATO",20,"LP1912"),
    @("07:15:48","07:43","215A_EL PATO",28,"LP1912"),
    @("07:40:11","07:44","215A_EL PATO",4,"LP1912"),
    @("08:40:59","08:43","215C_EL PATO",3,"LP1912"),
    @("08:30:14","08:44","215C_EL PATO",14,"LP1912"),
    @("08:52:33","08:53","215B_EL PATO",1,"LP1912"),
    @("08:30:14","08:54","215B_EL PATO",24,"LP1912"),
    @("08:40:59","08:57","215A_EL PATO",17,"LP1912"),
    @("08:52:33","08:58","215A_EL PATO",6,"LP1912"),
    @("09:23:52","09:58","215C_EL PATO",35,"LP1912"),
    @("08:52:33","09:59","215C_EL PATO",67,"LP1912"),
    @("10:41:48","11:30","215C_EL PATO",49,"LP1912"),
    @("11:13:01","11:31","215C_EL PATO",18,"LP1912"),
    @("11:35:40","11:41","215B_EL PATO",6,"LP1912"),
    @("11:13:01","11:42","215B_EL PATO",29,"LP1912"),
    @("11:35:40","12:21","215C_EL PATO",46,"LP1912"),
    @("12:12:04","12:22","215C_EL PATO",10,"LP1912"),
    @("12:54:41","13:03","215C_EL PATO",9,"LP1912"),
    @("12:47:00","13:04","215C_EL PATO",17,"LP1912"),
    @("12:33:54","13:32","215A_EL PATO",59,"LP1912"),
    @("13:14:41","13:33","215A_EL PATO",19,"LP1912"),
    @("14:33:43","14:33","215C_EL PATO",0,"LP1912"),
    @("14:15:21","14:34","215C_EL PATO",19,"LP1912"),
    @("14:33:43","14:47","215B_EL PATO",14,"LP1912"),
    @("14:47:05","14:48","215B_EL PATO",1,"LP1912"),
    @("14:53:58","14:53","215_EL PELIGRO",0,"LP1912"),
    @("14:47:05","14:54","215_EL PELIGRO",7,"LP1912"),
    @("15:17:56","15:33","215C_EL PATO",16,"LP1912"),
    @("14:47:05","15:34","215C_EL PATO",47,"LP1912"),
    @("15:59:49","17:09","215C_EL PATO",70,"LP1912"),
    @("15:59:49","17:38","215B_EL PATO",99,"LP1912"),
    @("15:47:47","17:42","215_EL PELIGRO",115,"LP1912"),
    @("15:59:49","17:45","215_EL PELIGRO",106,"LP1912")
)

for ($i = 0; $i -lt $data2.Length; $i++) {
    $row = 6 + $i
    $r = $data2[$i]
    $ws2.Cells.Item($row, 1).Value = $r[0]
    $ws2.Cells.Item($row, 2).Value = $r[1]
    $ws2.Cells.Item($row, 3).Value = $r[2]
    $ws2.Cells.Item($row, 4).Value = $r[3]
    $ws2.Cells.Item($row, 5).Value = $r[4]
}

# ---------------------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Cells.Item(2,1).Value = "Última actualización: $nuevaHora"

$data3 = @(
    @("07:15:48","07:42","215A_LA PLATA",27,"L6173"),
    @("07:40:11","07:43","215A_LA PLATA",3,"L6173"),
    @("07:15:48","08:35","215A_LA PLATA",80,"L6173"),
    @("08:30:14","08:36","215A_LA PLATA",6,"L6173"),
    @("08:40:59","08:50","215C_LA PLATA",10,"L6203"),
    @("08:30:14","08:51","215C_LA PLATA",21,"L6203"),
    @("08:52:33","08:52","215C_LA PLATA",0,"L6203"),
    @("08:40:59","09:20","215A_LA PLATA",40,"L6173"),
    @("08:52:33","09:21","215A_LA PLATA",29,"L6173"),
    @("09:23:52","09:23","215A_LA PLATA",0,"L6173"),
    @("10:07:51","10:12","215C_LA PLATA",5,"L6203"),
    @("08:52:33","10:13","215C_LA PLATA",81,"L6203"),
    @("10:07:51","10:29","215B_LP-P MOR-1 Y 57",22,"L6173"),
    @("08:52:33","10:30","215B_LP-P MOR-1 Y 57",98,"L6173"),
    @("10:07:51","10:30","215A_LA PLATA",23,"L6173"),
    @("08:52:33","10:31","215A_LA PLATA",99,"L6173"),
    @("10:41:48","11:25","215C_LA PLATA",44,"L6203"),
    @("11:13:01","11:26","215C_LA PLATA",13,"L6203"),
    @("12:54:41","13:11","215C_LA PLATA",17,"L6203"),
    @("12:47:00","13:12","215C_LA PLATA",25,"L6203"),
    @("13:14:41","13:16","215C_LA PLATA",2,"L6203"),
    @("13:14:41","13:20","215B_LP-P MOR-1 Y 57",6,"L6173"),
    @("12:47:00","13:21","215B_LP-P MOR-1 Y 57",34,"L6173"),
    @("13:14:41","13:56","215C_LA PLATA",42,"L6203"),
    @("13:43:25","13:57","215C_LA PLATA",14,"L6203"),
    @("13:57:31","13:58","215C_LA PLATA",1,"L6203"),
    @("14:15:21","14:27","215C_LA PLATA",12,"L6203"),
    @("14:33:43","14:34","215C_LA PLATA",1,"L6203"),
    @("15:17:56","15:21","215A_LA PLATA",4,"L6173"),
    @("14:47:05","15:22","215A_LA PLATA",35,"L6173"),
    @("15:59:49","16:01","215C_LA PLATA",2,"L6203"),
    @("14:47:05","16:02","215C_LA PLATA",75,"L6203"),
    @("15:59:49","16:29","215B_LP-P MOR-40 Y 115",30,"L6173"),
    @("14:47:05","16:30","215B_LP-P MOR-40 Y 115",103,"L6173"),
    @("15:59:49","17:05","215C_LA PLATA",66,"L6203")
)

for ($i = 0; $i -lt $data3.Length; $i++) {
    $row = 6 + $i
    $r = $data3[$i]
    $ws3.Cells.Item($row, 1).Value = $r[0]
    $ws3.Cells.Item($row, 2).Value = $r[1]
    $ws3.Cells.Item($row, 3).Value = $r[2]
    $ws3.Cells.Item($row, 4).Value = $r[3]
    $ws3.Cells.Item($row, 5).Value = $r[4]
}

Write-Host "Update complete"
